$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 1052.9
$ws.Range("I99").Value = 361.2857
$ws.Range("J99").Value = 2666.6667
$ws.Range("K99").Value = 1083.8571
$ws.Range("L99").Value = 8000.000100000001
$ws.Range("M99").Value = 414.1428999999998
$ws.Range("N99").Value = -10996.0001
$ws.Range("H132").Value = 1588.5667
$ws.Range("I132").Value = 1503.4584
$ws.Range("J132").Value = 1929
$ws.Range("K132").Value = 4510.3752
$ws.Range("L132").Value = 5787
$ws.Range("M132").Value = -1980.3752
$ws.Range("N132").Value = -10847
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H138").Value = 2449.818
$ws.Range("J138").Value = 3034.2
$ws.Range("L138").Value = 9102.599999999999
$ws.Range("N138").Value = -19382.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2623.5293
$ws.Range("I2").Value = 1995.0769
$ws.Range("K2").Value = 1995.0769
$ws.Range("M2").Value = -1882.0769
$ws.Range("H32").Value = 4938.1377
$ws.Range("I32").Value = 4938.1377
$ws.Range("K32").Value = 4938.1377
$ws.Range("M32").Value = -4651.1377
$ws.Range("H45").Value = 3911.3333
$ws.Range("I45").Value = 1598.2
$ws.Range("J45").Value = 5563.5713
$ws.Range("K45").Value = 1598.2
$ws.Range("L45").Value = 5563.5713
$ws.Range("M45").Value = -1221.2
$ws.Range("N45").Value = -6317.5713
$ws.Range("H74").Value = 1449.9143
$ws.Range("I74").Value = 1438.591
$ws.Range("J74").Value = 1469.0769
$ws.Range("K74").Value = 1438.591
$ws.Range("L74").Value = 1469.0769
$ws.Range("M74").Value = -564.5909999999999
$ws.Range("N74").Value = -3217.0769
$ws.Range("H77").Value = 1449.9143
$ws.Range("I77").Value = 1438.591
$ws.Range("J77").Value = 1469.0769
$ws.Range("K77").Value = 7192.955
$ws.Range("L77").Value = 7345.3845
$ws.Range("M77").Value = -2824.955
$ws.Range("N77").Value = -16081.3845
$ws.Range("H97").Value = 2333.3333
$ws.Range("I97").Value = 2032.0588
$ws.Range("J97").Value = 3613.75
$ws.Range("K97").Value = 2032.0588
$ws.Range("L97").Value = 3613.75
$ws.Range("M97").Value = -1536.0588
$ws.Range("N97").Value = -4605.75
$ws.Range("H116").Value = 2623.5293
$ws.Range("I116").Value = 1995.0769
$ws.Range("K116").Value = 1995.0769
$ws.Range("M116").Value = 298.9231
$ws.Range("H132").Value = 1501.6333
$ws.Range("I132").Value = 1333.1904
$ws.Range("J132").Value = 1894.6666
$ws.Range("K132").Value = 3999.5712
$ws.Range("L132").Value = 5683.9998
$ws.Range("M132").Value = -1469.5712
$ws.Range("N132").Value = -10743.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2623.5293
$ws.Range("I3").Value = 1995.0769
$ws.Range("K3").Value = 1995.0769
$ws.Range("M3").Value = -1881.0769
$ws.Range("H99").Value = 1120.7
$ws.Range("I99").Value = 1256
$ws.Range("J99").Value = 354
$ws.Range("K99").Value = 1256
$ws.Range("L99").Value = 354
$ws.Range("M99").Value = 242
$ws.Range("N99").Value = -3350
$ws.Range("H134").Value = 2469.5806
$ws.Range("I134").Value = 2702.8518
$ws.Range("J134").Value = 895
$ws.Range("K134").Value = 8108.555399999999
$ws.Range("L134").Value = 2685
$ws.Range("M134").Value = -5573.555399999999
$ws.Range("N134").Value = -7755

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2765.25
$ws.Range("I31").Value = 2899
$ws.Range("J31").Value = 2631.5
$ws.Range("K31").Value = 2899
$ws.Range("L31").Value = 2631.5
$ws.Range("M31").Value = -2604
$ws.Range("N31").Value = -3221.5
$ws.Range("H34").Value = 2765.25
$ws.Range("I34").Value = 2899
$ws.Range("J34").Value = 2631.5
$ws.Range("K34").Value = 2899
$ws.Range("L34").Value = 2631.5
$ws.Range("M34").Value = -2697
$ws.Range("N34").Value = -3035.5
$ws.Range("H107").Value = 1968.0769
$ws.Range("I107").Value = 676.8570999999999
$ws.Range("K107").Value = 676.8570999999999
$ws.Range("M107").Value = 1243.1429
$ws.Range("H132").Value = 1329.3334
$ws.Range("I132").Value = 1140.7142
$ws.Range("J132").Value = 1989.5
$ws.Range("K132").Value = 3422.1426
$ws.Range("L132").Value = 5968.5
$ws.Range("M132").Value = -892.1425999999997
$ws.Range("N132").Value = -11028.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 3622.5
$ws.Range("I56").Value = 3622.5
$ws.Range("K56").Value = 3622.5
$ws.Range("M56").Value = -3092.5
$ws.Range("H133").Value = 5000
$ws.Range("I133").Value = 5000
$ws.Range("K133").Value = 15000
$ws.Range("M133").Value = -9940

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2105.3215
$ws.Range("I132").Value = 2095.05
$ws.Range("K132").Value = 6285.150000000001
$ws.Range("M132").Value = -3755.150000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1004268
$ws.Range("I100").Value = 1501400
$ws.Range("J100").Value = 10004
$ws.Range("K100").Value = 1501400
$ws.Range("L100").Value = 10004
$ws.Range("M100").Value = -1500859
$ws.Range("N100").Value = -11086

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 3450
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H7").Value = 2005
$ws.Range("I7").Value = 2005
$ws.Range("K7").Value = 2005
$ws.Range("M7").Value = -1892
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("H11").Value = 16350
$ws.Range("J11").Value = 28200
$ws.Range("L11").Value = 28200
$ws.Range("N11").Value = -28484
$ws.Range("H13").Value = 9298.666999999999
$ws.Range("I13").Value = 1450
$ws.Range("K13").Value = 1450
$ws.Range("M13").Value = -1310
$ws.Range("H17").Value = 9941
$ws.Range("I17").Value = 9941
$ws.Range("K17").Value = 9941
$ws.Range("M17").Value = -9769
$ws.Range("H30").Value = 24383.2
$ws.Range("I30").Value = 22000
$ws.Range("J30").Value = 24979
$ws.Range("K30").Value = 22000
$ws.Range("L30").Value = 24979
$ws.Range("M30").Value = -21893
$ws.Range("N30").Value = -25193
